$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 3832.8333
$ws.Range("J92").Value = 6966.6665
$ws.Range("L92").Value = 6966.6665
$ws.Range("N92").Value = -9462.666499999999

# Hunk 1: sheet ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 333341000
$ws.Range("I98").Value = 500006240
$ws.Range("J98").Value = 10500
$ws.Range("K98").Value = 500006240
$ws.Range("L98").Value = 10500
$ws.Range("M98").Value = -500004742
$ws.Range("N98").Value = -13496

# Hunk 2: sheet ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 333341000
$ws.Range("I122").Value = 500006240
$ws.Range("J122").Value = 10500
$ws.Range("K122").Value = 1500018720
$ws.Range("L122").Value = 31500
$ws.Range("M122").Value = -1500016270
$ws.Range("N122").Value = -36400

# Hunk 3: sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1990.4237
$ws.Range("I132").Value = 1954.2142
$ws.Range("J132").Value = 2666.3333
$ws.Range("K132").Value = 5862.642599999999
$ws.Range("L132").Value = 7998.999899999999
$ws.Range("M132").Value = -3332.642599999999
$ws.Range("N132").Value = -13058.9999

# Hunk 4: sheet ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 715086.3
$ws.Range("I135").Value = 833771.5600000001
$ws.Range("K135").Value = 7503944.040000001
$ws.Range("M135").Value = -7501409.040000001

# Hunk 5: sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2954.7727
$ws.Range("I137").Value = 3535.5557
$ws.Range("K137").Value = 10606.6671
$ws.Range("M137").Value = -8056.667099999999

# Hunk 6: sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5539.6
$ws.Range("I138").Value = 2264.4707
$ws.Range("J138").Value = 7004.7896
$ws.Range("K138").Value = 6793.4121
$ws.Range("L138").Value = 21014.3688
$ws.Range("M138").Value = -1653.4121
$ws.Range("N138").Value = -31294.3688

# Hunk 7: sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1669413.9
$ws.Range("I32").Value = 1814232.6
$ws.Range("K32").Value = 1814232.6
$ws.Range("M32").Value = -1813945.6

# Hunk 8: sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5088.4443
$ws.Range("I61").Value = 1351.3448
$ws.Range("K61").Value = 1351.3448
$ws.Range("M61").Value = -1139.3448

# Hunk 9: sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4571.431
$ws.Range("I132").Value = 1938
$ws.Range("J132").Value = 9399.388999999999
$ws.Range("K132").Value = 5814
$ws.Range("L132").Value = 28198.167
$ws.Range("M132").Value = -3284
$ws.Range("N132").Value = -33258.167

# Hunk 10: sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5088.4443
$ws.Range("I136").Value = 1351.3448
$ws.Range("K136").Value = 4054.0344
$ws.Range("M136").Value = -1504.0344

# Hunk 11: sheet BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8334478.5
$ws.Range("I20").Value = 12821729
$ws.Range("K20").Value = 12821729
$ws.Range("M20").Value = -12821482

# Hunk 12: sheet BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6994361
$ws.Range("I99").Value = 862
$ws.Range("K99").Value = 862
$ws.Range("M99").Value = 636

# Hunk 13: sheet BSM row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 41547.668
$ws.Range("J100").Value = 41547.668
$ws.Range("L100").Value = 41547.668
$ws.Range("N100").Value = -43711.668

# Hunk 14: sheet BSM row 106
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 86773.664
$ws.Range("J106").Value = 86773.664
$ws.Range("L106").Value = 86773.664
$ws.Range("N106").Value = -89297.664

# Hunk 15: sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10144.363
$ws.Range("I134").Value = 12241.154
$ws.Range("J134").Value = 8781.450000000001
$ws.Range("K134").Value = 36723.462
$ws.Range("L134").Value = 26344.35
$ws.Range("M134").Value = -34188.462
$ws.Range("N134").Value = -31414.35

# Hunk 16: sheet CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2273470.2
$ws.Range("J12").Value = 3572034.5
$ws.Range("L12").Value = 10716103.5
$ws.Range("N12").Value = -10716449.5

# Hunk 17: sheet CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 9324.6
$ws.Range("I99").Value = 2541
$ws.Range("K99").Value = 7623
$ws.Range("M99").Value = -5377

# Hunk 18: sheet GSM row 29
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1253.5
$ws.Range("I29").Value = 1007
$ws.Range("J29").Value = 1500
$ws.Range("K29").Value = 1007
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = -717
$ws.Range("N29").Value = -2080

# Hunk 19: sheet GSM row 101
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 49574
$ws.Range("J101").Value = 49574
$ws.Range("L101").Value = 49574
$ws.Range("N101").Value = -56064

# Hunk 20: sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3244.6
$ws.Range("I102").Value = 3244.6
$ws.Range("K102").Value = 3244.6
$ws.Range("M102").Value = -1622.6

# Hunk 21: sheet GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8727.272000000001
$ws.Range("J113").Value = 9277.777
$ws.Range("L113").Value = 9277.777
$ws.Range("N113").Value = -13617.777

# Hunk 22: sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 14005.944
$ws.Range("J22").Value = 25888.875
$ws.Range("L22").Value = 25888.875
$ws.Range("N22").Value = -26478.875

# Hunk 23: sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 14005.944
$ws.Range("J27").Value = 25888.875
$ws.Range("L27").Value = 25888.875
$ws.Range("N27").Value = -26102.875

# Hunk 24: sheet LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3433.3333
$ws.Range("I32").Value = 3250
$ws.Range("J32").Value = 3800
$ws.Range("K32").Value = 3250
$ws.Range("L32").Value = 3800
$ws.Range("M32").Value = -2933
$ws.Range("N32").Value = -4434

# Hunk 25: sheet LTW row 38
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Hunk 26: sheet LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 58823880
$ws.Range("J55").Value = 516.8182
$ws.Range("L55").Value = 516.8182
$ws.Range("N55").Value = -862.8182

# Hunk 27: sheet LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 100002050
$ws.Range("I68").Value = 200000770
$ws.Range("J68").Value = 3339.6
$ws.Range("K68").Value = 200000770
$ws.Range("L68").Value = 3339.6
$ws.Range("M68").Value = -200000021
$ws.Range("N68").Value = -4837.6

# Hunk 28: sheet LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 100002050
$ws.Range("I71").Value = 200000770
$ws.Range("J71").Value = 3339.6
$ws.Range("K71").Value = 1000003850
$ws.Range("L71").Value = 16698
$ws.Range("M71").Value = -1000000106
$ws.Range("N71").Value = -24186

# Hunk 29: sheet LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4518.467
$ws.Range("I93").Value = 5129.6665
$ws.Range("J93").Value = 3601.6667
$ws.Range("K93").Value = 5129.6665
$ws.Range("L93").Value = 3601.6667
$ws.Range("M93").Value = -3881.6665
$ws.Range("N93").Value = -6097.6667

# Hunk 30: sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16137044
$ws.Range("I132").Value = 33338306
$ws.Range("K132").Value = 100014918
$ws.Range("M132").Value = -100012388

# Hunk 31: sheet WVR row 32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# Hunk 32: sheet WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 783.9286
$ws.Range("I107").Value = 588.7273
$ws.Range("K107").Value = 1766.1819
$ws.Range("M107").Value = 153.8181

# Hunk 33: sheet WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1654.4138
$ws.Range("I113").Value = 1479.5
$ws.Range("K113").Value = 4438.5
$ws.Range("M113").Value = -2268.5
